$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1 (existing title slide): merge the two subtitle runs into one.
# Force an actual text mutation first so the engine rebuilds the run
# (keeps the original hu-HU language + dirty/smtClean tracking), then
# set the final desired text.
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2).TextFrame.TextRange
$subtitle.Text = "Készítette: Beréti Zsófia."
$subtitle.Text = "Készítette: Beréti Zsófia"

$layoutText = [Microsoft.Office.Interop.PowerPoint.PpSlideLayout]::ppLayoutText
$layoutTitleOnly = [Microsoft.Office.Interop.PowerPoint.PpSlideLayout]::ppLayoutTitleOnly
$alignCenter = [Microsoft.Office.Interop.PowerPoint.PpParagraphAlignment]::ppAlignCenter
$mouseClick = [Microsoft.Office.Interop.PowerPoint.PpMouseActivation]::ppMouseClick

# ---------------------------------------------------------------------
# Slide 2 - "Cyber mezőgazdaság"
# ---------------------------------------------------------------------
$s2 = $p.Slides.Add(2, $layoutText)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Cyber mezőgazdaság"

# ---------------------------------------------------------------------
# Slide 3 - "jövő"
# ---------------------------------------------------------------------
$s3 = $p.Slides.Add(3, $layoutText)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "jövő"

# ---------------------------------------------------------------------
# Slide 4 - "előnyök"
# ---------------------------------------------------------------------
$s4 = $p.Slides.Add(4, $layoutText)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "előnyök"

# ---------------------------------------------------------------------
# Slide 5 - "Cloud (felhő) technológia" (bold)
# ---------------------------------------------------------------------
$s5 = $p.Slides.Add(5, $layoutText)
$t5 = $s5.Shapes.Item(1).TextFrame.TextRange
$t5.Text = "Cloud (felhő) technológia"
$t5.Font.Bold = [Microsoft.Office.Interop.PowerPoint.MsoTriState]::msoTrue

# ---------------------------------------------------------------------
# Slide 6 - "precíziós gazdálkodás" (bold)
# ---------------------------------------------------------------------
$s6 = $p.Slides.Add(6, $layoutText)
$t6 = $s6.Shapes.Item(1).TextFrame.TextRange
$t6.Text = "precíziós gazdálkodás"
$t6.Font.Bold = [Microsoft.Office.Interop.PowerPoint.MsoTriState]::msoTrue

# ---------------------------------------------------------------------
# Slide 7 - "GPS rendszerek" (bold)
# ---------------------------------------------------------------------
$s7 = $p.Slides.Add(7, $layoutText)
$t7 = $s7.Shapes.Item(1).TextFrame.TextRange
$t7.Text = "GPS rendszerek"
$t7.Font.Bold = [Microsoft.Office.Interop.PowerPoint.MsoTriState]::msoTrue

# ---------------------------------------------------------------------
# Slide 8 - "távérzékelés" (bold)
# ---------------------------------------------------------------------
$s8 = $p.Slides.Add(8, $layoutText)
$t8 = $s8.Shapes.Item(1).TextFrame.TextRange
$t8.Text = "távérzékelés"
$t8.Font.Bold = [Microsoft.Office.Interop.PowerPoint.MsoTriState]::msoTrue

# ---------------------------------------------------------------------
# Slides 9, 10, 11 - blank title+content slides
# ---------------------------------------------------------------------
$s9 = $p.Slides.Add(9, $layoutText)
$s10 = $p.Slides.Add(10, $layoutText)
$s11 = $p.Slides.Add(11, $layoutText)

# ---------------------------------------------------------------------
# Slide 12 - "Köszönöm a figyelmet!" (title-only layout, moved + centered)
# ---------------------------------------------------------------------
$s12 = $p.Slides.Add(12, $layoutTitleOnly)
$sh12 = $s12.Shapes.Item(1)
$sh12.Left = 825500 / 12700.0
$sh12.Top = 2333625 / 12700.0
$sh12.Width = 10515600 / 12700.0
$sh12.Height = 1325563 / 12700.0
$tr12 = $sh12.TextFrame.TextRange
$tr12.Text = "Köszönöm a figyelmet!"
$tr12.ParagraphFormat.Alignment = $alignCenter

# ---------------------------------------------------------------------
# Slide 13 - "Források:" with hyperlinked source list
# ---------------------------------------------------------------------
$s13 = $p.Slides.Add(13, $layoutText)
$s13.Shapes.Item(1).TextFrame.TextRange.Text = "Források:"

$urls = @(
    "https://mfor.hu/cikkek/befektetes/az-informatika-a-mezogazdasag-jovoje-nagyon-sok-mulhat-rajta.html",
    "https://infobex.hu/hirek/cyber-mezogazdasag-az-it-terhoditasa-az-agrikulturaban/hu",
    "https://agrarium7.hu/cikkek/2064-informatika-az-agrariumban",
    "https://innoskart.digital/innoskart_magazin/informatika-az-allattenyesztesben-agrarmernok-az-informatikusok-kozott/",
    "https://www.youtube.com/watch?v=1evSfdmUw34",
    "https://www.youtube.com/watch?v=xMt2mgYFwVc"
)

$body13 = $s13.Shapes.Item(2).TextFrame.TextRange
$body13.Text = "$($urls[0])`r$($urls[1])`r$($urls[2])`r$($urls[3])`r$($urls[4])`r$($urls[5])`r"

for ($i = 1; $i -le 5; $i++) {
    $para = $body13.Paragraphs($i, 1)
    $para.ActionSettings.Item($mouseClick).Hyperlink.Address = $urls[$i - 1]
}
